$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a value that Excel would otherwise auto-parse as a number
# (e.g. a bare "62%") to be stored as literal text, while preserving the
# cell's existing style (border/format). We enter it as a text formula
# ( ="62%" ) which Excel evaluates to the plain string, then convert that
# formula result to a static value in place via copy / paste-special-values.
function Set-CellText($cellRef, $text) {
    $ws.Range($cellRef).Formula = '="' + $text + '"'
    $ws.Range($cellRef).Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
}
$excel.CutCopyMode = $false

$ws.Range('E2').Value = '2026-02-23 21:48:29'
$ws.Range('E3').Value = '2026-02-23 21:48:31'
$ws.Range('E4').Value = '2026-02-23 21:48:34'
$ws.Range('O4').Value = '12.0 °C'
$ws.Range('E5').Value = '2026-02-23 21:48:36'
$ws.Range('O5').Value = '4.8 °C'
$ws.Range('E6').Value = '2026-02-23 21:48:38'
Set-CellText 'H6' '62%'
$ws.Range('O6').Value = '14.0 °C'
$ws.Range('E7').Value = '2026-02-23 21:48:40'
$ws.Range('E8').Value = '2026-02-23 21:48:42'
$ws.Range('E9').Value = '2026-02-23 21:48:45'
$ws.Range('O9').Value = '12.4 °C'
$ws.Range('E10').Value = '2026-02-23 21:48:47'
$ws.Range('O10').Value = '10.7 °C'
$ws.Range('E11').Value = '2026-02-23 21:48:49'
$ws.Range('O11').Value = '8.8 °C'
$ws.Range('E12').Value = '2026-02-23 21:48:51'
$ws.Range('O12').Value = '10.9 °C'
$ws.Range('E13').Value = '2026-02-23 21:48:54'
$ws.Range('G13').Value = '3 cm'
$ws.Range('J13').Value = '1026.9 hPa'
$ws.Range('O13').Value = '7.0 °C'
$ws.Range('E14').Value = '2026-02-23 21:48:56'
Set-CellText 'H14' '76%'
$ws.Range('O14').Value = '12.4 °C'
$ws.Range('E15').Value = '2026-02-23 21:48:59'
Set-CellText 'H15' '71%'
$ws.Range('O15').Value = '12.4 °C'
$ws.Range('E16').Value = '2026-02-23 21:49:01'
$ws.Range('E17').Value = '2026-02-23 21:49:03'
Set-CellText 'H17' '43%'
$ws.Range('E18').Value = '2026-02-23 21:49:06'
Set-CellText 'H18' '74%'
$ws.Range('O18').Value = '10.9 °C'
$ws.Range('E19').Value = '2026-02-23 21:49:08'
Set-CellText 'H19' '47%'
$ws.Range('O19').Value = '12.3 °C'
$ws.Range('E20').Value = '2026-02-23 21:49:10'
$ws.Range('O20').Value = '4.1 °C'
$ws.Range('E21').Value = '2026-02-23 21:49:13'
$ws.Range('J21').Value = '1025.7 hPa'
$ws.Range('O21').Value = '9.7 °C'
$ws.Range('E22').Value = '2026-02-23 21:49:15'
$ws.Range('E23').Value = '2026-02-23 21:49:18'
$ws.Range('E24').Value = '2026-02-23 21:49:20'
Set-CellText 'H24' '83%'
$ws.Range('E25').Value = '2026-02-23 21:49:22'
Set-CellText 'H25' '28%'
$ws.Range('O25').Value = '5.8 °C'
$ws.Range('E26').Value = '2026-02-23 21:49:25'
Set-CellText 'H26' '52%'
$ws.Range('E27').Value = '2026-02-23 21:49:27'
$ws.Range('E28').Value = '2026-02-23 21:49:30'
$ws.Range('J28').Value = '1025.0 hPa'
$ws.Range('O28').Value = '11.0 °C'
$ws.Range('E29').Value = '2026-02-23 21:49:32'
$ws.Range('O29').Value = '10.7 °C'
$ws.Range('E30').Value = '2026-02-23 21:49:34'
$ws.Range('O30').Value = '12.9 °C'
$ws.Range('E31').Value = '2026-02-23 21:49:37'
$ws.Range('J31').Value = '1023.7 hPa'
$ws.Range('E32').Value = '2026-02-23 21:49:39'
Set-CellText 'H32' '68%'
$ws.Range('O32').Value = '7.4 °C'
$ws.Range('E33').Value = '2026-02-23 21:49:41'
$ws.Range('J33').Value = '1025.3 hPa'
$ws.Range('O33').Value = '8.7 °C'
$ws.Range('E34').Value = '2026-02-23 21:49:43'
Set-CellText 'H34' '43%'
$ws.Range('O34').Value = '4.0 °C'
$ws.Range('E35').Value = '2026-02-23 21:49:46'
$ws.Range('N35').Value = '6.7 °C 21:10 TU'
$ws.Range('O35').Value = '12.1 °C'
$ws.Range('E36').Value = '2026-02-23 21:49:48'
Set-CellText 'H36' '73%'
$ws.Range('O36').Value = '13.0 °C'
$ws.Range('E37').Value = '2026-02-23 21:49:51'
Set-CellText 'H37' '67%'
$ws.Range('O37').Value = '9.0 °C'
$ws.Range('E38').Value = '2026-02-23 21:49:53'
Set-CellText 'H38' '65%'
$ws.Range('O38').Value = '12.2 °C'
$ws.Range('E39').Value = '2026-02-23 21:49:56'
Set-CellText 'H39' '26%'
$ws.Range('E40').Value = '2026-02-23 21:49:58'
$ws.Range('O40').Value = '8.7 °C'
$ws.Range('E41').Value = '2026-02-23 21:50:00'
$ws.Range('O41').Value = '11.9 °C'
$ws.Range('E42').Value = '2026-02-23 21:50:03'
$ws.Range('E43').Value = '2026-02-23 21:50:05'
$ws.Range('E44').Value = '2026-02-23 21:50:08'
Set-CellText 'H44' '35%'
$ws.Range('E45').Value = '2026-02-23 21:50:10'
$ws.Range('O45').Value = '8.2 °C'
$ws.Range('E46').Value = '2026-02-23 21:50:12'
$ws.Range('O46').Value = '10.2 °C'

$excel.CutCopyMode = $false
